# Update the cached "datetimeFigureOut" field text (10/19/2024 -> 10/20/2024)
# on every Date placeholder across the slide master and all of its custom
# layouts, then fix the typo in the "Currently?" shape title on slide 1.

$p = $ppt.ActivePresentation

$oldDate = "10/19/2024"
$newDate = "10/20/2024"

function Update-DateShape($shp) {
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master's own Date placeholder.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

# Every custom layout under the master has its own Date placeholder too.
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}

# Slide 1: drop the stray "?" from the "Currently?" box title.
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "Currently?: LLM and GNN Security and Applications") {
                $shp.TextFrame.TextRange.Text = "Currently: LLM and GNN Security and Applications"
            }
        }
    }
}
